$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.765.18'
$ws.Range('E2').Value = '  +1.74%  '
$ws.Range('D3').Value = '1.879.39'
$ws.Range('E3').Value = '  +1.34%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = "'332.83"
$ws.Range('D6').Value = "'1.004"
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').Value = "'0.4715"
$ws.Range('E7').Value = '  +3.93%  '
$ws.Range('E8').Value = '  +2.38%  '
$ws.Range('D9').Value = "'47.88"
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('D10').Value = "'0.08060"
$ws.Range('E10').Value = '  +1.79%  '
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('D12').Value = "'22.25"
$ws.Range('E12').Value = '  +4.16%  '
$ws.Range('D13').Value = '1.873.59'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = "'5.984"
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').Value = "'7.128"
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').Value = "'1.004"
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('E17').Value = '  +1.69%  '
$ws.Range('E18').Value = '  +2.29%  '
$ws.Range('D19').Value = "'0.06675"
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('D20').Value = "'17.31"
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = '27.785.60'
$ws.Range('E22').Value = '  +1.79%  '
$ws.Range('D23').Value = "'5.535"
$ws.Range('E23').Value = '  +0.64%  '
$ws.Range('D24').Value = "'11.01"
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('D25').Value = "'2.311"
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').Value = '2.099.64'
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('D27').Value = "'159.47"
$ws.Range('E27').Value = '  +3.79%  '
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('E29').Value = '  +2.40%  '
$ws.Range('E30').Value = '  +2.17%  '
$ws.Range('D31').Value = "'121.95"
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('E32').Value = '  +5.19%  '
$ws.Range('D33').Value = "'0.09534"
$ws.Range('E33').Value = '  +2.35%  '
$ws.Range('D34').Value = "'1.449"
$ws.Range('E34').Value = '  -0.64%  '
$ws.Range('D35').Value = "'3.598"
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('D36').Value = "'5.362"
$ws.Range('E36').Value = '  +1.83%  '
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('E38').Value = '  +1.62%  '
$ws.Range('D39').Value = "'1.235"
$ws.Range('E39').Value = '  +1.23%  '
$ws.Range('D40').Value = "'8.140"
$ws.Range('E40').Value = '  +1.07%  '
$ws.Range('D41').Value = "'0.6031"
$ws.Range('E41').Value = '  +2.09%  '
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').Value = "'0.1906"
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('D44').Value = "'10.27"
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('D45').Value = "'1.276"
$ws.Range('D46').Value = "'0.5716"
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('D47').Value = "'12.27"
$ws.Range('E47').Value = '  +1.85%  '
$ws.Range('D48').Value = "'1.951"
$ws.Range('E48').Value = '  +1.86%  '
$ws.Range('D49').Value = "'3.380"
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').Value = "'0.06924"
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('D51').Value = "'114.37"
$ws.Range('E51').Value = '  +5.27%  '
